$wb = $excel.ActiveWorkbook

# --- "Time Range" sheet updates ---
$ws = $wb.Worksheets.Item("Time Range")

# Header: rename numUniqueYears -> numUniqueTimeSeriesEntries
$ws.Range("E1").Value = "numUniqueTimeSeriesEntries"

# Updated counts (now counting unique time series entries instead of unique years)
$ws.Range("E2").Value = 48
$ws.Range("E3").Value = 348
$ws.Range("E5").Value = 335
$ws.Range("E6").Value = 746

# Column E needs to be wider to fit the new header text
# (stored OOXML width = ColumnWidth + 5/6, so this yields width="25")
$ws.Columns.Item(5).ColumnWidth = 24.166666666666668

# Update the remembered selection on this sheet
$ws.Range("I9").Select()

# --- "POD v POU" sheet updates ---
$ws1 = $wb.Worksheets.Item("POD v POU")
$ws1.Range("F14:F15").Select()
